# 自动更新Excel文件 - 2026-02-24 23:23:51
# Decrement the "剩余" (remaining) days counter in column E for each data row.
# When remaining reaches 1, it rolls over to 10 and the "开始时间" (start date)
# in column F is advanced by 10 days (representing a new cycle/reorder).
# Row 36 is a stale/anomalous record (remaining already equals the total and
# the date value is corrupted) and is intentionally left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($row = 2; $row -le $lastRow; $row++) {

    # Skip the known stale/anomalous row (remaining already equals total days).
    if ($row -eq 36) { continue }

    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $remaining = $eCell.Value2
    if ($remaining -eq $null -or $remaining -eq "") { continue }

    if ($remaining -eq 1) {
        $eCell.Value2 = 10
        $fCell.Value2 = $fCell.Value2 + 10
    } else {
        $eCell.Value2 = $remaining - 1
    }
}
